$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'25.918.47"
$ws.Cells.Item(2, 5).Value = "  -1.14%  "
$ws.Cells.Item(3, 4).Value = "'1.636.29"
$ws.Cells.Item(3, 5).Value = "  -0.52%  "
$ws.Cells.Item(4, 5).Value = "  +0.11%  "
$ws.Cells.Item(5, 4).Value = "'215.45"
$ws.Cells.Item(5, 5).Value = "  -0.61%  "
$ws.Cells.Item(6, 5).Value = "  +0.02%  "
$ws.Cells.Item(7, 5).Value = "  +0.08%  "
$ws.Cells.Item(8, 5).Value = "  -0.74%  "
$ws.Cells.Item(9, 4).Value = "'0.0637"
$ws.Cells.Item(9, 5).Value = "  -0.07%  "
$ws.Cells.Item(10, 5).Value = "  -1.39%  "
$ws.Cells.Item(11, 4).Value = "'0.0793"
$ws.Cells.Item(11, 5).Value = "  +0.18%  "
$ws.Cells.Item(12, 5).Value = "  -0.06%  "
$ws.Cells.Item(13, 4).Value = "'1.863.16"
$ws.Cells.Item(13, 5).Value = "  -0.54%  "
$ws.Cells.Item(14, 4).Value = "'1.615.89"
$ws.Cells.Item(14, 5).Value = "  -1.48%  "
$ws.Cells.Item(15, 5).Value = "  -0.27%  "
$ws.Cells.Item(16, 5).Value = "  -0.19%  "
$ws.Cells.Item(17, 4).Value = "'62.83"
$ws.Cells.Item(17, 5).Value = "  -0.68%  "
$ws.Cells.Item(18, 4).Value = "'25.924.77"
$ws.Cells.Item(18, 5).Value = "  -1.16%  "
$ws.Cells.Item(19, 5).Value = "  +0.18%  "
$ws.Cells.Item(20, 5).Value = "  -1.31%  "
$ws.Cells.Item(21, 5).Value = "  -1.74%  "
$ws.Cells.Item(22, 4).Value = "'9.92"
$ws.Cells.Item(22, 5).Value = "  -1.33%  "
$ws.Cells.Item(23, 5).Value = "  -0.43%  "
$ws.Cells.Item(24, 4).Value = "'1.81"
$ws.Cells.Item(24, 5).Value = "  +1.14%  "
$ws.Cells.Item(25, 5).Value = "  +4.74%  "
$ws.Cells.Item(26, 4).Value = "'144.08"
$ws.Cells.Item(26, 5).Value = "  +0.80%  "
$ws.Cells.Item(27, 5).Value = "  -0.10%  "
$ws.Cells.Item(28, 5).Value = "  -0.61%  "
$ws.Cells.Item(29, 4).Value = "'15.53"
$ws.Cells.Item(29, 5).Value = "  -0.73%  "
$ws.Cells.Item(30, 4).Value = "'1.24"
$ws.Cells.Item(30, 5).Value = "  -0.43%  "
$ws.Cells.Item(31, 4).Value = "'0.0500"
$ws.Cells.Item(31, 5).Value = "  -0.36%  "
$ws.Cells.Item(32, 5).Value = "  -2.01%  "
$ws.Cells.Item(33, 4).Value = "'3.24"
$ws.Cells.Item(33, 5).Value = "  -0.14%  "
$ws.Cells.Item(34, 5).Value = "  -3.83%  "
$ws.Cells.Item(35, 5).Value = "  +1.41%  "
$ws.Cells.Item(36, 5).Value = "  -1.11%  "
$ws.Cells.Item(37, 4).Value = "'1.132.72"
$ws.Cells.Item(37, 5).Value = "  -0.16%  "
$ws.Cells.Item(38, 5).Value = "  -1.58%  "
$ws.Cells.Item(39, 5).Value = "  -1.95%  "
$ws.Cells.Item(40, 5).Value = "  -0.37%  "
$ws.Cells.Item(41, 4).Value = "'5.50"
$ws.Cells.Item(41, 5).Value = "  +0.00%  "
$ws.Cells.Item(42, 5).Value = "  -0.88%  "
$ws.Cells.Item(43, 4).Value = "'0.794"
$ws.Cells.Item(43, 5).Value = "  -0.46%  "
$ws.Cells.Item(44, 4).Value = "'1.772.72"
$ws.Cells.Item(44, 5).Value = "  -0.65%  "
$ws.Cells.Item(45, 5).Value = "  +2.40%  "
$ws.Cells.Item(46, 4).Value = "'56.63"
$ws.Cells.Item(46, 5).Value = "  -0.63%  "
$ws.Cells.Item(47, 5).Value = "  +2.28%  "
$ws.Cells.Item(48, 4).Value = "'1.46"
$ws.Cells.Item(48, 5).Value = "  -0.76%  "
$ws.Cells.Item(49, 4).Value = "'7.67"
$ws.Cells.Item(49, 5).Value = "  -0.05%  "
$ws.Cells.Item(50, 5).Value = "  -0.87%  "
$ws.Cells.Item(51, 4).Value = "'0.0959"
$ws.Cells.Item(51, 5).Value = "  -1.07%  "
